# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4):
#   F2: 8613 -> 8640
#   F4: 395  -> 398
#   F5: 29   -> 30

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 8640
    $ws.Range("F4").Value = 398
    $ws.Range("F5").Value = 30
}
